$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# This status string appears on the Overview sheet (columns E/F, row 2)
# and on each per-locale sheet (column C, row 2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes ---
# Target raw (OOXML) column width is 13.4101845877511 characters.
# The Excel COM ColumnWidth property is offset from the raw/stored width
# by the standard 5-pixel padding expressed in characters (5/6 for this
# engine's column metrics), so we back that out before assigning.
$targetColumnWidth = 13.4101845877511 - (5 / 6)

# Overview sheet: columns E and F (zh-cn / de-de handoff-date columns)
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

# zh-cn sheet: column C (Status column)
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth

# de-de sheet: column C (Status column)
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
